$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression (only B2 changes slightly)
$ws.Range("B2").Value = 0.1749439332307385

# Row 3: RandomForestRegressor
$ws.Range("B3").Value = 0.02289963500604794
$ws.Range("C3").Value = 0.02263920663709551
$ws.Range("D3").Value = 0.04829514876820767

# Row 4: model renamed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.0244326463740279
$ws.Range("C4").Value = 0.02400890201203397
$ws.Range("D4").Value = 0.04160863179126655

# Row 5: model renamed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.0213756523163499
$ws.Range("C5").Value = 0.02386392024273786
$ws.Range("D5").Value = 0.02295920102333551
